$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.929.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.737.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.59%  '
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.382'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.223.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.761.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000149'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.744.72'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '353.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.519'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.13%  '
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0903'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.96'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.982'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '345.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.03'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0582'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.623'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0250'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.139.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.00%  '
